$wb = $excel.ActiveWorkbook

# --- Sheet 1: Cases by Age Group ---
$ws = $wb.Worksheets.Item("Cases by Age Group")
$ws.Range("B2").Value = 274
$ws.Range("B3").Value = 1352
$ws.Range("B4").Value = 3718
$ws.Range("B5").Value = 15636
$ws.Range("B6").Value = 17170
$ws.Range("B7").Value = 15078
$ws.Range("B8").Value = 12694
$ws.Range("B9").Value = 4582
$ws.Range("B10").Value = 3095
$ws.Range("B11").Value = 1875
$ws.Range("B12").Value = 1218
$ws.Range("B13").Value = 1909
$ws.Range("G9").Select()

# --- Sheet 2: Cases by Gender ---
$ws = $wb.Worksheets.Item("Cases by Gender")
$ws.Range("B2").Value = 26641
$ws.Range("B3").Value = 51075
$ws.Range("B4").Value = 898
$ws.Activate()
$ws.Range("E22").Select()
$excel.ActiveWindow.Zoom = 84

# --- Sheet 3: Cases by RaceEthnicity ---
$ws = $wb.Worksheets.Item("Cases by RaceEthnicity")
$ws.Range("B2").Value = 956
$ws.Range("B3").Value = 12965
$ws.Range("B4").Value = 28237
$ws.Range("B5").Value = 536
$ws.Range("B6").Value = 27305
$ws.Range("B7").Value = 8615
$ws.Range("B8").Select()

# --- Sheet 4: Fatalities by Age Group ---
$ws = $wb.Worksheets.Item("Fatalities by Age Group")
$ws.Range("B2").Value = 7
$ws.Range("B5").Value = 243
$ws.Range("B6").Value = 820
$ws.Range("B7").Value = 2401
$ws.Range("B8").Value = 5505
$ws.Range("B9").Value = 4605
$ws.Range("B10").Value = 5973
$ws.Range("B11").Value = 6599
$ws.Range("B12").Value = 6520
$ws.Range("B13").Value = 16499
$ws.Range("E5").Select()

# --- Sheet 5: Fatalities by Gender ---
$ws = $wb.Worksheets.Item("Fatalities by Gender")
$ws.Range("B2").Value = 20638
$ws.Range("B3").Value = 28578
$ws.Range("E6").Select()

# --- Sheet 6: Fatalities by Race-Ethnicity ---
$ws = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$ws.Range("B2").Value = 1037
$ws.Range("B3").Value = 4948
$ws.Range("B4").Value = 22868
$ws.Range("B5").Value = 269
$ws.Range("B6").Value = 20072
$ws.Range("B7").Value = 23
$ws.Range("B11").Select()

$wb.Worksheets.Item("Cases by Age Group").Activate()
